# Update countries & provincias Spain
# Applies the COVID data refresh captured in the commit: refreshed case
# counts for several countries and the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Septiembre de 2020 a las 21:38"

# --- Update numeric data cells (row -> column -> new value) ---
# Row 4: Estados Unidos
$ws.Range("B4").Value = 7025605
$ws.Range("C4").Value = 20837
$ws.Range("D4").Value = 4282651
$ws.Range("E4").Value = 2538659
$ws.Range("G4").Value = 177
$ws.Range("H4").Value = 204295

# Row 5: India
$ws.Range("B5").Value = 5557573
$ws.Range("C5").Value = 71961
$ws.Range("D5").Value = 4492145
$ws.Range("E5").Value = 976485

# Row 14: Francia
$ws.Range("D14").Value = 93008
$ws.Range("E14").Value = 333715

# Row 29: Canada
$ws.Range("B29").Value = 144686
$ws.Range("C29").Value = 1037
$ws.Range("D29").Value = 125215
$ws.Range("E29").Value = 10248
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = 9223

# Row 120: Hong Kong
$ws.Range("B120").Value = 5045
$ws.Range("C120").Value = 266
$ws.Range("D120").Value = 3527
$ws.Range("E120").Value = 1486
$ws.Range("G120").Value = 2
$ws.Range("H120").Value = 32

# Row 121: Congo
$ws.Range("B121").Value = 5039
$ws.Range("C121").Value = 6
$ws.Range("D121").Value = 4717
$ws.Range("E121").Value = 219
$ws.Range("H121").Value = 103

# Row 122: Guinea Ecuatorial
$ws.Range("C122").Value = 16
$ws.Range("D122").Value = 3887
$ws.Range("E122").Value = 1026
$ws.Range("H122").Value = 89

# Row 123: Nicaragua
$ws.Range("B123").Value = 5002
$ws.Range("D123").Value = 4509
$ws.Range("E123").Value = 410
$ws.Range("H123").Value = 83

# Row 124: Republica de Africa Central
$ws.Range("B124").Value = 4961
$ws.Range("D124").Value = 2913
$ws.Range("E124").Value = 1901
$ws.Range("H124").Value = 147

# Row 125: Jordania
$ws.Range("B125").Value = 4786
$ws.Range("D125").Value = 1830
$ws.Range("E125").Value = 2894
$ws.Range("H125").Value = 62

# Row 131: Lituania
$ws.Range("B131").Value = 3833
$ws.Range("C131").Value = 33
$ws.Range("D131").Value = 963
$ws.Range("E131").Value = 2695
$ws.Range("G131").Value = 3
$ws.Range("H131").Value = 175

# Row 132: Siria
$ws.Range("B132").Value = 3814
$ws.Range("C132").Value = 70
$ws.Range("D132").Value = 2199
$ws.Range("E132").Value = 1528
$ws.Range("H132").Value = 87

# Row 153: Yemen
$ws.Range("B153").Value = 2028
$ws.Range("C153").Value = 2
$ws.Range("D153").Value = 1235
$ws.Range("E153").Value = 207

# Row 214: Montserrat
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

# Row 215: Islas Malvinas
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
